# Contest 21: PBKS vs KKR - enter player scores for row 30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E30").Value = 0
$ws.Range("H30").Value = 60
$ws.Range("K30").Value = 80
$ws.Range("N30").Value = 40
$ws.Range("Q30").Value = 20
$ws.Range("T30").Value = 100
